$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '95.924.10'
$ws.Range('E2').Value = '  +0.52%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.563.86'
$ws.Range('E3').Value = '  -1.83%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.91'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '654.03'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.61'
$ws.Range('E7').Value = '  +10.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.405'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.05'
$ws.Range('E9').Value = '  +5.42%  '
$ws.Range('B10').Value = 'USDC'
$ws.Range('C10').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.00'
$ws.Range('E10').Value = '  +0.06%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '3.561.74'
$ws.Range('E11').Value = '  -1.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.22'
$ws.Range('E12').Value = '  +1.03%  '
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.37'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.225.24'
$ws.Range('E15').Value = '  -2.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '95.815.36'
$ws.Range('E16').Value = '  +0.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000258'
$ws.Range('E17').Value = '  +1.55%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.565.55'
$ws.Range('E18').Value = '  -1.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.61'
$ws.Range('E20').Value = '  -3.35%  '
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.515'
$ws.Range('E22').Value = '  +6.52%  '
$ws.Range('B23').Value = 'SuiNetwork'
$ws.Range('C23').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.38'
$ws.Range('E23').Value = '  -6.22%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '501.96'
$ws.Range('E24').Value = '  -1.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.91'
$ws.Range('E25').Value = '  +4.33%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000197'
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '95.71'
$ws.Range('E27').Value = '  -0.48%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.78'
$ws.Range('E28').Value = '  +0.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.755.55'
$ws.Range('E29').Value = '  -1.52%  '
$ws.Range('E30').Value = '  +9.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.00'
$ws.Range('E31').Value = '  -4.93%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.29'
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.181'
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.17%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '31.24'
$ws.Range('E36').Value = '  -3.41%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '8.76'
$ws.Range('E37').Value = '  +7.32%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '612.04'
$ws.Range('E38').Value = '  +6.76%  '
$ws.Range('E39').Value = '  -0.01%  '
$ws.Range('E40').Value = '  +9.38%  '
$ws.Range('E42').Value = '  -0.13%  '
$ws.Range('E43').Value = '  -3.87%  '
$ws.Range('E44').Value = '  +4.64%  '
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.53'
$ws.Range('E46').Value = '  -0.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0420'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.26'
$ws.Range('E48').Value = '  +0.73%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '33.64'
$ws.Range('E49').Value = '  -5.36%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.51'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.16'
$ws.Range('E51').Value = '  +0.95%  '
